# Apply updated ocean-freight surcharge values (hybrid/electric transport
# calculator increases) to the OceanFreightDatabase worksheet, and leave the
# active selection where the author last left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new price value (column D = PRICE)
$updates = @{
    "D3"  = 1200
    "D4"  = 1300
    "D8"  = 1200
    "D9"  = 1300
    "D13" = 1200
    "D14" = 1300
    "D18" = 1500
    "D19" = 1600
    "D20" = 1700
    "D23" = 1500
    "D24" = 1600
    "D25" = 1700
    "D42" = 1700
    "D47" = 1700
    "D52" = 2000
    "D53" = 2100
    "D54" = 2200
    "D56" = 1000
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value2 = $updates[$cell]
}

# Match the author's final cursor/selection position on the sheet.
$ws.Activate()
$ws.Range("D62").Select()
